$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "settings" sheet
#    - BOT_NAME value: "PushyReminder" -> "EndlessReminder"
#    - new row inserted for CHATWORK_API_TOKEN (between SLACK_ICON_EMOJI
#      and TIME_INTERVAL)
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("settings")

$wsSettings.Range("B2").Value = "EndlessReminder"

$wsSettings.Rows("6").Insert()
$wsSettings.Range("A6").Value = "CHATWORK_API_TOKEN"

# ---------------------------------------------------------------------
# 2. "main" sheet
#    - duplicate row 2 into a new row 3 to add a second schedule entry
#      that reminds by day of the week ("Fri") instead of day of month
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("main")

$wsMain.Rows("2").Copy()
$wsMain.Rows("3").Insert()

$wsMain.Range("A3").Value = 2
$wsMain.Range("D3").Value = "Fri"
$wsMain.Range("E3").Value = $false

# ---------------------------------------------------------------------
# 3. Restore view/selection state seen in the final workbook:
#    settings sheet selection sits on A10, main sheet stays the active
#    tab with A3 selected and scrolled back to the top-left.
# ---------------------------------------------------------------------
$wsSettings.Range("A10").Select() | Out-Null

$wsMain.Activate() | Out-Null
$wsMain.Range("A1").Select() | Out-Null
$wsMain.Range("A3").Select() | Out-Null
